# Tasks.xlsx update: add a new task row to the table and turn on the
# table's Total Row (Sum of Hours), matching the "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- 1. Fill in the (already table-reserved) row 7 with the new task ---
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "spsavebaneditlog deadlog analyst"
$ws.Cells.Item(7, 3).Value = 42992
$ws.Cells.Item(7, 4).Value = 6
$ws.Cells.Item(7, 5).Value = "Vista"

# Give the new "When" cell the same date style/format as the rest of the column.
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(7, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Turn on the table's Total Row ---
$tbl.ShowTotals = $true

# Totals row "When" cell keeps the date styling (empty, like Excel leaves it).
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(8, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Totals row "Hours" cell: Sum formula.
$ws.Cells.Item(8, 4).Formula = "=SUM(D2:D7)"

# --- 3. Match the trailing selection left behind in the worksheet ---
$ws.Range("F17").Select() | Out-Null
